$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.963.54'
$ws.Range('E2').Value = '  -0.96%  '

$ws.Range('D3').Value = '2.336.94'
$ws.Range('E3').Value = '  +1.20%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.94'
$ws.Range('E5').Value = '  -1.36%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.50'
$ws.Range('E6').Value = '  -1.67%  '

$ws.Range('E7').Value = '  -4.30%  '

$ws.Range('E8').Value = '  +0.07%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.509'
$ws.Range('E9').Value = '  -3.78%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.78'
$ws.Range('E10').Value = '  -4.55%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.42'
$ws.Range('E11').Value = '  +1.12%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0796'
$ws.Range('E12').Value = '  -2.22%  '

$ws.Range('E13').Value = '  +0.77%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.84'
$ws.Range('E14').Value = '  -2.93%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.78'
$ws.Range('E15').Value = '  +4.71%  '

$ws.Range('D16').Value = '2.316.45'
$ws.Range('E16').Value = '  +0.26%  '

$ws.Range('E17').Value = '  +2.18%  '

$ws.Range('D18').Value = '42.901.27'
$ws.Range('E18').Value = '  -0.87%  '

$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0911'
$ws.Range('E19').Value = '  -2.24%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.73'
$ws.Range('E20').Value = '  -4.89%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.18'
$ws.Range('E21').Value = '  -0.19%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.22'
$ws.Range('E22').Value = '  +0.12%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.47'
$ws.Range('E23').Value = '  -2.15%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.01'
$ws.Range('E24').Value = '  -0.56%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.55'
$ws.Range('E25').Value = '  -2.31%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.05%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.60'
$ws.Range('E27').Value = '  +3.17%  '

$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.32'
$ws.Range('E28').Value = '  +2.48%  '

$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '35.45'
$ws.Range('E29').Value = '  -3.69%  '

$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.28'
$ws.Range('E30').Value = '  -3.91%  '

$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '161.63'
$ws.Range('E31').Value = '  -4.05%  '

$ws.Range('B32').Value = 'FirstDigitalUSD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  -0.03%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.11'
$ws.Range('E33').Value = '  -3.47%  '

$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.64'
$ws.Range('E34').Value = '  +4.23%  '

$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.47'
$ws.Range('E35').Value = '  -2.91%  '

$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.44'
$ws.Range('E36').Value = '  -3.69%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0723'
$ws.Range('E37').Value = '  -2.66%  '

$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.85'
$ws.Range('E38').Value = '  -1.88%  '

$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.90'
$ws.Range('E39').Value = '  -4.98%  '

$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.102'
$ws.Range('E40').Value = '  -3.49%  '

$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.113'
$ws.Range('E41').Value = '  -2.87%  '

$ws.Range('B42').Value = 'ApeXProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.59'
$ws.Range('E42').Value = '  +7.00%  '

$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.021.97'
$ws.Range('E43').Value = '  +2.06%  '

$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0284'
$ws.Range('E44').Value = '  -3.76%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.87'
$ws.Range('E45').Value = '  -1.98%  '

$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.18'
$ws.Range('E46').Value = '  +1.69%  '

$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.93'
$ws.Range('E47').Value = '  -2.82%  '

$ws.Range('B48').Value = 'MultiversX'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.92'
$ws.Range('E48').Value = '  +0.01%  '

$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.89'
$ws.Range('E49').Value = '  -1.97%  '

$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.562.02'
$ws.Range('E50').Value = '  +1.12%  '

$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.71'
$ws.Range('E51').Value = '  +0.82%  '
